$d = $word.ActiveDocument

$d.Paragraphs(8).Range.Text  = "Design: Yes"
$d.Paragraphs(12).Range.Text = "Design: No problems"
$d.Paragraphs(16).Range.Text = "Design: Budget was good, yard was refunded."
$d.Paragraphs(20).Range.Text = "Design: First time in 20 years automation design made all drawings, not outsourced. "
$d.Paragraphs(24).Range.Text = "Design: Machinery design have skilled designers."
$d.Paragraphs(28).Range.Text = "Design: Jira and ERM  doesn't work as I expected. ACAD without electrical symbols."
$d.Paragraphs(32).Range.Text = "Design: More teamwork before purchase of different sister systems. "
